$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P mirrors the style of column O for the header and data rows
$ws.Range("O4:O5").Copy()
$ws.Range("P4:P5").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("P4").Value = 2021
$ws.Range("P5").Value = 80.9

$ws.Range("N10").Select()
